$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing rows 2-4 and add new rows 5-8 (permeation frames data)

# Row 2: ion_id 2333
$ws.Range("A2").NumberFormat = "@"
$ws.Range("A2").Value = "2333"
$ws.Range("A2").Style = "Normal"
$ws.Range("B2").Value = 1248
$ws.Range("C2").Value = "[-3.3383673429489136, 11.95009446144104, -6.190477780997753]"
$ws.Range("D2").Value = 13.86619880539935
$ws.Range("E2").Value = 9.180969760765738
$ws.Range("F2").Value = 0.6621115050788658
$ws.Range("G2").Value = 12.40763692867548
$ws.Range("H2").Value = -6.190477780997753
$ws.Range("I2").Value = "[-1.2147483825683594, 0.6887054443359375, -0.09078216552734375]"

# Row 3: ion_id 2343
$ws.Range("A3").NumberFormat = "@"
$ws.Range("A3").Value = "2343"
$ws.Range("A3").Style = "Normal"
$ws.Range("B3").Value = 1186
$ws.Range("C3").Value = "[2.4428126215934753, 0.49578909622505307, -8.41463577747345]"
$ws.Range("D3").Value = 8.7760603689633
$ws.Range("E3").Value = 8.280676746972738
$ws.Range("F3").Value = 0.9435528470448431
$ws.Range("G3").Value = 2.492617165180414
$ws.Range("H3").Value = -8.41463577747345
$ws.Range("I3").Value = "[4.815853118896484, -3.690704345703125, -12.438545227050781]"

# Row 4: ion_id 2372
$ws.Range("A4").NumberFormat = "@"
$ws.Range("A4").Value = "2372"
$ws.Range("A4").Style = "Normal"
$ws.Range("B4").Value = 1248
$ws.Range("C4").Value = "[-9.981672696769238, -4.027546465396881, -5.427530646324158]"
$ws.Range("D4").Value = 12.05458457488881
$ws.Range("E4").Value = -10.895716383016
$ws.Range("F4").Value = -0.9038649416184051
$ws.Range("G4").Value = 10.76359235368747
$ws.Range("H4").Value = -5.427530646324158
$ws.Range("I4").Value = "[3.4089317321777344, 1.604522705078125, 0.10662841796875]"

# Row 5: ion_id 2334
$ws.Range("A5").NumberFormat = "@"
$ws.Range("A5").Value = "2334"
$ws.Range("A5").Style = "Normal"
$ws.Range("B5").Value = 1248
$ws.Range("C5").Value = "[3.699047952890396, -14.163208246231079, -8.058695323765278]"
$ws.Range("D5").Value = 16.70993698088557
$ws.Range("E5").Value = -2.613813113001731
$ws.Range("F5").Value = -0.1564226792711224
$ws.Range("G5").Value = 14.63828622427812
$ws.Range("H5").Value = -8.058695323765278
$ws.Range("I5").Value = "[-5.621417999267578, 1.4537582397460938, -3.01348876953125]"

# Row 6: ion_id 2230
$ws.Range("A6").NumberFormat = "@"
$ws.Range("A6").Value = "2230"
$ws.Range("A6").Style = "Normal"
$ws.Range("B6").Value = 1225
$ws.Range("C6").Value = "[-0.9318812191486359, -2.220571478828788, -13.284280061721802]"
$ws.Range("D6").Value = 13.50079394174466
$ws.Range("E6").Value = 13.21380134611099
$ws.Range("F6").Value = 0.978742539374201
$ws.Range("G6").Value = 2.408181948937791
$ws.Range("H6").Value = -13.2842800617218
$ws.Range("I6").Value = "[1.012481689453125, -0.3558006286621094, -11.524864196777344]"

# Row 7: ion_id 2515
$ws.Range("A7").NumberFormat = "@"
$ws.Range("A7").Value = "2515"
$ws.Range("A7").Style = "Normal"
$ws.Range("B7").Value = 1219
$ws.Range("C7").Value = "[-7.849294036626816, 10.417173475027084, -12.752289831638336]"
$ws.Range("D7").Value = 18.24143130437794
$ws.Range("E7").Value = 5.720300019439035
$ws.Range("F7").Value = 0.3135883321867491
$ws.Range("G7").Value = 13.04334773293395
$ws.Range("H7").Value = -12.75228983163834
$ws.Range("I7").Value = "[9.918651580810547, -4.4640045166015625, -19.940650939941406]"

# Row 8: ion_id 2280
$ws.Range("A8").NumberFormat = "@"
$ws.Range("A8").Value = "2280"
$ws.Range("A8").Style = "Normal"
$ws.Range("B8").Value = 1248
$ws.Range("C8").Value = "[13.266657590866089, 2.8086227253079414, -3.7235074639320374]"
$ws.Range("D8").Value = 14.06261259796197
$ws.Range("E8").Value = 0.5318098541636465
$ws.Range("F8").Value = 0.03781728682767804
$ws.Range("G8").Value = 13.56069929046438
$ws.Range("H8").Value = -3.723507463932037
$ws.Range("I8").Value = "[-0.21262359619140625, 0.6027069091796875, -0.4115142822265625]"
